$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data appended to the bottom of the sheet (rows 207-209)
$data = @(
    @(46043, 2110.53, 114302.64, 112192.11, 6875.6201171875, 1.172786593437195),
    @(46044, 2110.53, 114865.89, 112755.36, 6913.35009765625, 1.167297005653381),
    @(46045, 2110.53, 113153.61, 111043.08, 6915.60986328125, 1.175461053848267)
)

$startRow = 207
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}

# Carry the formatting of the last pre-existing row (206, style "s=2" on
# column A i.e. the date column) down onto the freshly appended rows so the
# new cells keep the same cell style as the rest of the table.
$ws.Range("A206").Copy() | Out-Null
$ws.Range("A207:A209").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
